# Author's commit: ":speech_balloon: Update by feedback"
#
# The diff removes the two trailing retrospective slides that used to
# sit right before the closing "Thank you" slide (old slide positions
# 14 and 15), leaving the "Thank you" slide as the new slide 14. (The
# diff's cached Notes Master datetimeFigureOut text bump is just
# PowerPoint re-caching that field on save and isn't independently
# settable through the object model, so it's left alone here.)

$p = $ppt.ActivePresentation

# Delete the higher-indexed slide first so the lower index stays valid.
$p.Slides.Item(15).Delete()
$p.Slides.Item(14).Delete()
